$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Smoke_Suite")

# The Incident ID column (C) test-data values were refreshed to newly
# created ServiceNow incident numbers as part of the Problem Management
# test-data changes.
$ws.Range("C2").Value = "INC0021759"
$ws.Range("C3").Value = "INC0021760"
$ws.Range("C4").Value = "INC0021759"
$ws.Range("C5").Value = "INC0021761"
